$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row: CAJEROS / 444 / 123
# A4 keeps the sheet's default (no explicit) style, just like A1/A4 in the
# source, so set it directly.
$ws.Range("A4").Value = "CAJEROS"

# B4/C4 use the same "text" cell style (number format 49 + border) already
# used by B2/C2, so copy that formatting across before writing the values.
# Writing the values AFTER the format is applied keeps "444"/"123" stored as
# text (shared strings) instead of being auto-converted to numbers.
$ws.Range("B2:C2").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B4").Value = "444"
$ws.Range("C4").Value = "123"

# Move the active selection to A3 (was C3 before the edit)
$ws.Range("A3").Select()

# The workbook was also resized in Excel; reflect the new window size.
$excel.ActiveWindow.Width = 14025
$excel.ActiveWindow.Height = 4710
